# Edit: add "BAU" to the Involved_Scenarios column on Uncertainty_Table,
# and add a new "Region" column (= "UGA") to the Setup sheet; change the
# number of parallel runs (D2) from 10 to 5. Mirrors the commit that
# prepares the workbook to copy results into the 1_Experiment folder.

$wb = $excel.ActiveWorkbook

# --- Uncertainty_Table: Involved_Scenarios (column D) now lists BAU too ---
$ws1 = $wb.Worksheets.Item("Uncertainty_Table")
for ($r = 2; $r -le 15; $r++) {
    $ws1.Cells.Item($r, 4).Value = "BAU ; Scenario1"
}
$ws1.Select()
$ws1.Range("D2:D15").Select()

# --- Setup: add a Region column, and change the parallel-run count ---
$ws2 = $wb.Worksheets.Item("Setup")
$ws2.Range("D2").Value = 5
$ws2.Range("J1").Value = "Region"
$ws2.Range("J2").Value = "UGA"
